$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Locoto" at Vega Modelo
# de Temuco. It sorts ahead of the row that was previously first (row 14),
# so insert a fresh row there; Excel shifts every row below it (old 14-38)
# down by one (new 15-39), preserving all of their data untouched.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new observation.
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44797
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100112042
$ws.Range("G14").Value = "Locoto"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 2700
$ws.Range("L14").Value = 2700
$ws.Range("M14").Value = 2700
$ws.Range("N14").Value = "$/kilo"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 2700
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
